# Apply updated submission/pass counts to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 1358
$ws.Range("F2").Value = 737

$ws.Range("E3").Value = 1634

$ws.Range("E4").Value = 605
$ws.Range("F4").Value = 226

$ws.Range("E5").Value = 423
$ws.Range("F5").Value = 201

$ws.Range("E9").Value = 513

$ws.Range("E11").Value = 401
$ws.Range("F11").Value = 155

$ws.Range("E13").Value = 304
$ws.Range("F13").Value = 204

$ws.Range("E14").Value = 273

$ws.Range("E19").Value = 116
$ws.Range("F19").Value = 66

$ws.Range("E50").Value = 35
$ws.Range("F50").Value = 5
